$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 2 (the "learningJob" parsing-job definition row) into row 3
# for the new "learningJob2" parsing job.

# First, paste formats only so every cell in the destination row gets
# materialised (avoids the "skip blank cells" behaviour of a values-only
# paste for the handful of source cells that hold an empty string).
$ws.Range("A2:CK2").Copy()
$ws.Range("A3:CK3").PasteSpecial(-4122)

# Now paste the real values/types (keeps the original text typing instead
# of Excel re-interpreting numeric-looking text like "0" / "100" as
# numbers).
$ws.Range("A2:CK2").Copy()
$ws.Range("A3:CK3").PasteSpecial(-4104)

$excel.CutCopyMode = $false

# Give the new job row its own unique name.
$ws.Range("A3").Value = "learningJob2"
